$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 80

# Use a scratch cell with a text formula to produce the plain text
# "01-07-2021", then copy/paste its value into A80. This avoids Excel's
# automatic date detection turning the literal string into a date
# serial number (which would otherwise alter the cell's number format).
$scratch = "Z1"
$ws.Range($scratch).Formula = '="01-07-2021"'
$ws.Range($scratch).Copy()
$ws.Range("A$row").PasteSpecial(-4163)
$ws.Range($scratch).Clear()

$ws.Range("B$row").Value = 39.2
$ws.Range("C$row").Value = 40.1
$ws.Range("D$row").Value = 45.5
$ws.Range("E$row").Value = 29.1
$ws.Range("F$row").Value = 37.2
$ws.Range("G$row").Value = 38.2
